$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 45: was the last ("open ended", no bottom border) row of the table.
# It becomes the last row of its group, so it gains a bottom border:
#   A45 (was empty, unstyled)         -> empty, style of "bottom border" row (like A7/A35/A37)
#   B45/C45/D45/E45 (no-border style) -> bottom-border style, same values
# ---------------------------------------------------------------------------
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A45:E45").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# New row 46 - single-line group (top+bottom border style, like rows 40/43)
# ---------------------------------------------------------------------------
$ws.Range("C46").Value = " Good luck catching that criminal\n[CS:N]Grovyle[CR]!"
$ws.Range("A46").Value = "SCRIPT/G01P03A/um1306.ssb"
$ws.Range("B46").Value = 490
$ws.Range("D46").Value = " Удачи вам в поимке этого\nнегодяя [CS:N]Гровайла[CR]!"
$ws.Range("E46").Value = " Ôäàœé âàí â ðïéíëå üóïãï\nîåãïäÿÿ [CS:N]Ãñïâàêìà[CR]!"

$ws.Range("A40:E40").Copy() | Out-Null
$ws.Range("A46:E46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(46).RowHeight = 43.2

# ---------------------------------------------------------------------------
# New row 47 - single-line group (top+bottom border style, like rows 40/43)
# ---------------------------------------------------------------------------
$ws.Range("C47").Value = " I hope you do well on your\n[CS:P]Crystal Cave[CR] exploration!"
$ws.Range("A47").Value = "SCRIPT/G01P03A/um1407.ssb"
$ws.Range("B47").Value = 450
$ws.Range("D47").Value = " Надеюсь, что у вас получится\nуспешно исследовать [CS:P]Кристальную Пещеру[CR]!"
$ws.Range("E47").Value = " Îàäåýòû, œóï ô âàò ðïìôœéóòÿ\nôòðåšîï éòòìåäïâàóû [CS:P]Ëñéòóàìûîôý Ðåþåñô[CR]!"

$ws.Range("A40:E40").Copy() | Out-Null
$ws.Range("A47:E47").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(47).RowHeight = 43.2

# ---------------------------------------------------------------------------
# New row 48 - final, open-ended row of the table (no-border style, like row 44)
# ---------------------------------------------------------------------------
$ws.Range("C48").Value = " It\'s back to our regular work\nfor now."
$ws.Range("A48").Value = "SCRIPT/G01P03A/um1603.ssb"
$ws.Range("B48").Value = 385
$ws.Range("D48").Value = " Пока что мы возвращаемся к нашей\nобычной работе."
$ws.Range("E48").Value = " Ðïëà œóï íú âïèâñàþàåíòÿ ë îàšåê\nïáúœîïê ñàáïóå."

$ws.Range("A44:E44").Copy() | Out-Null
$ws.Range("A48:E48").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(48).RowHeight = 47.4

# ---------------------------------------------------------------------------
# View state: scrolled down a bit further, selection moved to C48
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("C48").Select() | Out-Null
